$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B127").Value = "32.922276,35.081794"
$ws.Range("B145").Value = "32.139558,34.959151"
$ws.Range("B146").Value = "32.027141,34.898820"
$ws.Range("B156").Value = "32.780066,35.517715"
$ws.Range("B164").Value = "32.801454,35.069185"
$ws.Range("B170").Value = "32.780066,35.517715"
$ws.Range("B176").Value = "32.819327,34.997607"
$ws.Range("B183").Value = "32.603829,35.298516"
$ws.Range("B184").Value = "32.807619,35.057422"
$ws.Range("B188").Value = "32.782155,34.976622"
$ws.Range("B195").Value = "32.815143,35.060720"
$ws.Range("B199").Value = "32.471755,34.969755"
$ws.Range("B205").Value = "32.163217,34.961133"
$ws.Range("B206").Value = "32.163217,34.961133"
$ws.Range("B209").Value = "32.790077,35.516279"
$ws.Range("B211").Value = "32.139558,34.959151"
$ws.Range("B224").Value = "32.139558,34.959151"
$ws.Range("B268").Value = "32.139558,34.959151"
$ws.Range("B301").Value = "31.750585,35.215673"
$ws.Range("B354").Value = "32.163217,34.961133"
$ws.Range("B362").Value = "32.163217,34.961133"
$ws.Range("B372").Value = "31.755957,34.989832"
$ws.Range("B384").Value = "32.139558,34.959151"
$ws.Range("B387").Value = "32.139558,34.959151"
$ws.Range("B388").Value = "32.165553,34.813406"
$ws.Range("B390").Value = "31.982527,34.765084"
$ws.Range("B392").Value = "31.807623,34.664804"
$ws.Range("B398").Value = "32.055436,34.805472"
$ws.Range("B400").Value = "32.020682,34.805150"
$ws.Range("B401").Value = "32.097022,34.829235"
$ws.Range("B403").Value = "31.665784,34.601137"
$ws.Range("B407").Value = "31.677567,34.596921"
$ws.Range("B410").Value = "31.946849,34.879864"
$ws.Range("B414").Value = "31.858601,35.215336"
$ws.Range("B415").Value = "31.858484,35.215449"
$ws.Range("B417").Value = "31.750492,35.215772"
$ws.Range("B419").Value = "31.858484,35.215449"
$ws.Range("B420").Value = "31.792463,35.144323"
$ws.Range("B422").Value = "31.753295,34.996429"
$ws.Range("B424").Value = "32.033552,34.851439"
$ws.Range("B427").Value = "31.225747,34.809580"
$ws.Range("B434").Value = "31.862441,35.220615"
$ws.Range("B435").Value = "31.225747,34.809580"
$ws.Range("B441").Value = "31.223100,34.820208"
$ws.Range("B445").Value = "32.171208,34.826985"
$ws.Range("B452").Value = "31.226551,34.807177"
$ws.Range("B453").Value = "31.826014,34.658552"
$ws.Range("B455").Value = "31.666206,34.591622"
$ws.Range("B457").Value = "32.001232,34.801778"
$ws.Range("B458").Value = "32.095724,34.858840"
$ws.Range("B459").Value = "31.223512,34.880824"
$ws.Range("B460").Value = "31.217089,34.816739"
$ws.Range("B464").Value = "31.236692,34.796056"
$ws.Range("B465").Value = "32.068716,34.778964"
$ws.Range("B474").Value = "32.107402,34.938858"
$ws.Range("B475").Value = "31.665784,34.601137"
$ws.Range("B476").Value = "32.139558,34.959151"
$ws.Range("B483").Value = "32.047811,34.882122"

Write-Output "Updated coordinate values in column B."
